$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B2" = 11.33622468515689
    "C2" = 4.246804205994255
    "D2" = 9.042552701705244
    "E2" = 13.63193449268408
    "F2" = 33.75237425238953
    "J2" = 9.96955608028912
    "K2" = 10.71800481330014
    "N2" = 19.32907680124848
    "O2" = 25.54615279963864
    "B3" = 11.0833433801261
    "C3" = 4.023729007910307
    "D3" = 8.98200666143981
    "E3" = 13.56835971976751
    "F3" = 33.79094773599699
    "J3" = 9.97570660182655
    "K3" = 10.54684682131501
    "N3" = 19.39007248772496
    "O3" = 25.61521826754155
    "B4" = 10.92732458837231
    "C4" = 3.879456251352308
    "D4" = 8.946288665052165
    "E4" = 13.53206170138501
    "F4" = 33.82232951032452
    "J4" = 9.981048082998392
    "K4" = 10.44218786508059
    "N4" = 19.42926453858477
    "O4" = 25.66271270910524
    "B5" = 10.86365212720371
    "C5" = 3.818859765322129
    "D5" = 8.93211242421083
    "E5" = 13.51797013938435
    "F5" = 33.83705123336895
    "J5" = 9.983618710367256
    "K5" = 10.39970067960938
    "N5" = 19.44567458088835
    "O5" = 25.68334431743754
    "B6" = 10.85307654011682
    "C6" = 3.808689879186618
    "D6" = 8.929781726051914
    "E6" = 13.51567286425033
    "F6" = 33.83961245800074
    "J6" = 9.984069362787988
    "K6" = 10.3926570754595
    "N6" = 19.44842600683022
    "O6" = 25.68684725624409
    "B7" = 10.9264661306821
    "C7" = 3.87864628317807
    "D7" = 8.946095927977094
    "E7" = 13.53186880764533
    "F7" = 33.8225202278557
    "J7" = 9.981081155932534
    "K7" = 10.44161414063493
    "N7" = 19.42948407107269
    "O7" = 25.6629857856886
    "B8" = 11.24924244445613
    "C8" = 4.171421182082246
    "D8" = 9.021381674150639
    "E8" = 13.6094532337912
    "F8" = 33.76407565702539
    "J8" = 9.971352220041329
    "K8" = 10.65893162965606
    "N8" = 19.34974767665808
    "O8" = 25.56890929775288
    "B9" = 11.87227593528181
    "C9" = 4.686488946343938
    "D9" = 9.180000155021679
    "E9" = 13.78274926925017
    "F9" = 33.71061283972599
    "J9" = 9.964671247500929
    "K9" = 11.08620502099705
    "N9" = 19.20713362415489
    "O9" = 25.42488916265317
    "B10" = 12.31890842837118
    "C10" = 5.027773412406021
    "D10" = 9.302420473099005
    "E10" = 13.92215153614517
    "F10" = 33.70866542213526
    "J10" = 9.967288165549547
    "K10" = 11.39776038798471
    "N10" = 19.11064990517342
    "O10" = 25.34386567196772
    "B11" = 12.5187034924389
    "C11" = 5.174812138324579
    "D11" = 9.359208833258037
    "E11" = 13.98800305192882
    "F11" = 33.71588117223116
    "J11" = 9.970103536249225
    "K11" = 11.53836425012385
    "N11" = 19.06853957198291
    "O11" = 25.31241129402895
    "B12" = 12.59380008679443
    "C12" = 5.229301613752002
    "D12" = 9.380855407790515
    "E12" = 14.01327268560491
    "F12" = 33.71977668883922
    "J12" = 9.971402296710027
    "K12" = 11.59139645286065
    "N12" = 19.05284816355173
    "O12" = 25.30127896705365
    "B13" = 12.57765280612792
    "C13" = 5.217619425923322
    "D13" = 9.376187368525866
    "E13" = 14.00781588503468
    "F13" = 33.71888603194299
    "J13" = 9.971112256236177
    "K13" = 11.57998519228782
    "N13" = 19.05621627288575
    "O13" = 25.30364185155896
    "B14" = 12.52489339975292
    "C14" = 5.179318935079358
    "D14" = 9.36098695495814
    "E14" = 13.99007542536082
    "F14" = 33.71617836328576
    "J14" = 9.970205729879549
    "K14" = 11.54273175622748
    "N14" = 19.06724352933519
    "O14" = 25.31147981395896
    "B15" = 12.49250147512438
    "C15" = 5.155703439499116
    "D15" = 9.351694286106795
    "E15" = 13.97925173678235
    "F15" = 33.7146712315926
    "J15" = 9.969680720564352
    "K15" = 11.51988397577946
    "N15" = 19.07403119768929
    "O15" = 25.31638225841669
    "B16" = 12.3057767658625
    "C16" = 5.017997913542508
    "D16" = 9.298729971463954
    "E16" = 13.91789555211591
    "F16" = 33.70835673581789
    "J16" = 9.967136787591798
    "K16" = 11.38854467088096
    "N16" = 19.11343764292369
    "O16" = 25.34603017615566
    "B17" = 12.19030611210318
    "C17" = 4.931408578331538
    "D17" = 9.266508217621096
    "E17" = 13.88086750237039
    "F17" = 33.70655691845121
    "J17" = 9.965991661571204
    "K17" = 11.30764825626559
    "N17" = 19.13806742604474
    "O17" = 25.36560349759318
    "B18" = 12.1235764490747
    "C18" = 4.880832015915844
    "D18" = 9.248079505330658
    "E18" = 13.85980037940443
    "F18" = 33.70628444649606
    "J18" = 9.96548602966266
    "K18" = 11.26101581099817
    "N18" = 19.15240149405546
    "O18" = 25.37737007599648
    "B19" = 12.10093138632482
    "C19" = 4.863575297050631
    "D19" = 9.241858271195454
    "E19" = 13.85270750794798
    "F19" = 33.70632323756517
    "J19" = 9.965341140808762
    "K19" = 11.245210731808
    "N19" = 19.15728359670769
    "O19" = 25.3814413206601
    "B20" = 12.20263122074554
    "C20" = 4.940706208366094
    "D20" = 9.269927589393689
    "E20" = 13.8847854761342
    "F20" = 33.70666958308986
    "J20" = 9.966097732788375
    "K20" = 11.31627084152734
    "N20" = 19.13542819917445
    "O20" = 25.36346724090044
    "B21" = 12.54040592846632
    "C21" = 5.190601104266028
    "D21" = 9.365447957532929
    "E21" = 13.99527732908541
    "F21" = 33.71694212560449
    "J21" = 9.970465693599166
    "K21" = 11.55368011765265
    "N21" = 19.06399764936349
    "O21" = 25.30915646612126
    "B22" = 12.75785740905402
    "C22" = 5.346977516693893
    "D22" = 9.428696535029399
    "E22" = 14.06942325355131
    "F22" = 33.73043356641664
    "J22" = 9.974675854710672
    "K22" = 11.7075868806756
    "N22" = 19.01879863684462
    "O22" = 25.27820113454484
    "B23" = 12.64212528444798
    "C23" = 5.264154510522153
    "D23" = 9.394869814903331
    "E23" = 14.02967915956223
    "F23" = 33.7226136318071
    "J23" = 9.972305158607783
    "K23" = 11.62557456533171
    "N23" = 19.04278671384802
    "O23" = 25.29430665295654
    "B24" = 12.19706010495568
    "C24" = 4.936505225955912
    "D24" = 9.268381391175119
    "E24" = 13.88301347114573
    "F24" = 33.70661627284475
    "J24" = 9.966049302221322
    "K24" = 11.31237295320042
    "N24" = 19.13662085129473
    "O24" = 25.36443144278459
    "B25" = 11.70533160283831
    "C25" = 4.553606133814425
    "D25" = 9.135995941296491
    "E25" = 13.73368429218426
    "F25" = 33.718520901396
    "J25" = 9.965154271169791
    "K25" = 10.97081581670394
    "N25" = 19.24425161500128
    "O25" = 25.45950547494625
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
